$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1 "ROW50-FE-LIFTER": append new row 31 (after existing row 30)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$r = 31
$ws1.Cells.Item($r, 1).Value = 45739.650805
$ws1.Cells.Item($r, 1).NumberFormat = $ws1.Cells.Item($r - 1, 1).NumberFormat
$ws1.Cells.Item($r, 2).Value = "0x01,0x90"
$ws1.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
$ws1.Cells.Item($r, 4).Value = "0x01,0x76"
$ws1.Cells.Item($r, 5).Value = "0xe"
$ws1.Cells.Item($r, 6).Value = 400
$ws1.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws1.Cells.Item($r, 8).Value = 374
$ws1.Cells.Item($r, 9).Value = 14

# ---------------------------------------------------------------------
# Sheet 2 "ROW50-MID-LIFTER": append new row 33 (after existing row 32)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$r = 33
$ws2.Cells.Item($r, 1).Value = 45739.62341435185
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item($r - 1, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = "0x01,0x90 "
$ws2.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws2.Cells.Item($r, 4).Value = "0x01,0x7a"
$ws2.Cells.Item($r, 5).Value = "0x19"
$ws2.Cells.Item($r, 6).Value = 400
# This value has more significant digits than Excel's 15-digit numeric
# precision, so (as with the rest of column G on this sheet) it is kept
# as literal text rather than a rounded number.
$ws2.Cells.Item($r, 7).Value = "'568631262647113771663628"
$ws2.Cells.Item($r, 7).Style = "Normal"
$ws2.Cells.Item($r, 8).Value = 378
$ws2.Cells.Item($r, 9).Value = 25

# ---------------------------------------------------------------------
# Sheet 3 "ROW11-FE-LIFTER": append new row 31 (after existing row 30)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$r = 31
$ws3.Cells.Item($r, 1).Value = 45739.67029701389
$ws3.Cells.Item($r, 1).NumberFormat = $ws3.Cells.Item($r - 1, 1).NumberFormat
$ws3.Cells.Item($r, 2).Value = "0x01,0x90"
$ws3.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws3.Cells.Item($r, 4).Value = "0x01,0x76"
$ws3.Cells.Item($r, 5).Value = "0x14"
$ws3.Cells.Item($r, 6).Value = 400
$ws3.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws3.Cells.Item($r, 8).Value = 374
$ws3.Cells.Item($r, 9).Value = 20

# ---------------------------------------------------------------------
# Sheet 4 "ROW11-MID-LIFTER": append new row 31 (after existing row 30)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$r = 31
$ws4.Cells.Item($r, 1).Value = 45739.81708040509
$ws4.Cells.Item($r, 1).NumberFormat = $ws4.Cells.Item($r - 1, 1).NumberFormat
$ws4.Cells.Item($r, 2).Value = "0x01,0x90"
$ws4.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws4.Cells.Item($r, 4).Value = "0x01,0x7e"
$ws4.Cells.Item($r, 5).Value = "0x19"
$ws4.Cells.Item($r, 6).Value = 400
$ws4.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws4.Cells.Item($r, 8).Value = 382
$ws4.Cells.Item($r, 9).Value = 25
